$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Generalforsamlingsdato" (column A) is obsolete and was removed from the
# template. Deleting the whole column shifts B:I left to A:H, so
# Udbetalingsdato/Identifikation/Navn/C-O/Adresse/Postnr./Land/Bruttoudbytte
# now occupy A1:H1 and the dangling shared string is dropped automatically.
$ws.Columns("A").Delete()

# The saved cursor position in the fixed template is A2.
$ws.Range("A2").Select()
